# Updates cryptos list data (prices and 1h volume %) per upstream diff.
# Some rows also swap coin name/link (rank reordering).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.663.73"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "2.342.70"
$ws.Range("E3").Value = "  +2.83%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.57"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.651"
$ws.Range("E6").Value = "  +2.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "67.00"
$ws.Range("E7").Value = "  +5.34%  "
$ws.Range("E9").Value = "  +2.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0972"
$ws.Range("E10").Value = "  -3.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.52"
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "27.12"
$ws.Range("E12").Value = "  +1.89%  "
$ws.Range("D13").Value = "2.689.61"
$ws.Range("E13").Value = "  +2.64%  "
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.60"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.26"
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.853"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").Value = "2.348.71"
$ws.Range("E18").Value = "  +2.78%  "
$ws.Range("D19").Value = "43.661.35"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").Value = "0.0₃0980"
$ws.Range("E20").Value = "  -2.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.34"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.27"
$ws.Range("E22").Value = "  +2.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.90"
$ws.Range("E23").Value = "  -1.27%  "
$ws.Range("B24").Value = "WEMIXToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.81"
$ws.Range("E24").Value = "  +13.93%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.44"
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.00"
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "175.07"
$ws.Range("E29").Value = "  +1.81%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.26"
$ws.Range("E30").Value = "  +6.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.46"
$ws.Range("E31").Value = "  +5.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.130"
$ws.Range("E32").Value = "  -6.84%  "
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.02"
$ws.Range("E34").Value = "  +4.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0690"
$ws.Range("E35").Value = "  -1.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.99"
$ws.Range("E36").Value = "  +1.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.50"
$ws.Range("E37").Value = "  +7.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.58"
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.63"
$ws.Range("E39").Value = "  -4.52%  "
$ws.Range("E40").Value = "  -2.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.06"
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "18.34"
$ws.Range("E43").Value = "  +3.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.18"
$ws.Range("E44").Value = "  +8.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.68"
$ws.Range("E45").Value = "  +1.29%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.20"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0956"
$ws.Range("E47").Value = "  -1.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.37"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").Value = "1.447.14"
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("E51").Value = "  -4.74%  "
